$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix case data: the "Case_table_dependency" column (I) for rows 3-7 still
# pointed at the old "ps_blob_001" case id; it should reference the renamed
# "btree_ps_blob_001" case instead.
$ws.Range("I3:I7").Value = "btree_ps_blob_001"

# Restore the active cell/selection left by the author's edit.
$ws.Range("I11").Select()
